$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Row 6: InputCSVFilePath keeps its name (A6 unchanged), but Value/Description
# switch from an absolute Windows path to a relative one, and a Description
# is added.
$ws.Range("B6").Value = "Data\Input\GitHubRepoURLInput.csv"
$ws.Range("C6").Value = "A relative path pointing to a place to look for the input CSV file before asking the user for input."

# Row 7: was GitCloneRootFilePath -> now GitCloneRepoATargetPath
$ws.Range("A7").Value = "GitCloneRepoATargetPath"
$ws.Range("B7").Value = "Data\Temp\RepoA"
$ws.Range("C7").Value = "A relative path pointing to a place to clone the first repo in each transaction to be compared."

# Row 8: was GitDiffOutputFilePath -> now GitCloneRepoBTargetPath
$ws.Range("A8").Value = "GitCloneRepoBTargetPath"
$ws.Range("B8").Value = "Data\Temp\RepoB"
$ws.Range("C8").Value = "A relative path pointing to a place to clone thesecond repo in each transaction to be compared."

# Row 9 is new: GitDiffOutputFolderPath (description is entered before the
# value so the shared-string table gets the same ordering as the authored
# workbook).
$ws.Range("A9").Value = "GitDiffOutputFolderPath"
$ws.Range("C9").Value = "A relative path pointing to a place to store the diff outputs for the repositories to be compared."
$ws.Range("B9").Value = "Data\Temp"

# These rows picked up a plain/default cell style (distinct xf record) in
# the authored workbook - touching NumberFormat forces the creation of a
# fresh style index without altering the visible formatting.
$ws.Range("A6:C9").NumberFormat = "General"

# Reflect the final selection/active cell left behind in the saved file.
$ws.Activate() | Out-Null
$ws.Range("B10").Select() | Out-Null
